# Auto-generated script to apply cryptos.xlsx diff via Excel COM interop
# Updates Price (D) and Volume(1h) (E) columns, and for a couple of rows the
# Coin name (B) / Link (C) columns, matching the upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "97.385.22"
$ws.Range("E2").Value = "  +1.62%  "
$ws.Range("D3").Value = "3.594.88"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.80"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +17.88%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "652.89"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.29%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.426"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.62%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.10%  "
$ws.Range("E10").Value = "  +2.78%  "
$ws.Range("D11").Value = "3.592.50"
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "44.89"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.41%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.204"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.44%  "
$ws.Range("E14").Value = "  +0.01%  "
$ws.Range("D15").Value = "4.263.31"
$ws.Range("E15").Value = "  +0.09%  "
$ws.Range("D16").Value = "97.302.33"
$ws.Range("E16").Value = "  +1.66%  "
$ws.Range("E17").Value = "  +2.07%  "
$ws.Range("D18").Value = "3.604.45"
$ws.Range("E18").Value = "  +0.47%  "
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.65"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "18.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.52%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.530"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.99%  "
$ws.Range("B23").Value = "SuiNetwork"
$ws.Range("C23").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.89%  "
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "518.10"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000206"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.96"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.42%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "103.13"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "13.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.183"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +25.33%  "
$ws.Range("E30").Value = "  -1.86%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.03"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.83%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.191"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.998"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "31.88"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.72"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.585"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.47%  "
$ws.Range("E38").Value = "  +2.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "617.36"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.93%  "
$ws.Range("E40").Value = "  +1.51%  "
$ws.Range("E41").Value = "  +2.60%  "
$ws.Range("B42").Value = "USDe"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.930"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.51%  "
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.454"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +38.15%  "
$ws.Range("E45").Value = "  +4.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0452"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.35"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.65"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.81%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.71"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "32.91"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.29"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.51%  "
